$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Duplicate "30_actual" into a new sheet "30_actual_hard", placed right
#    after "30_actual" (this is what Excel's Worksheet.Copy does when the
#    "after" sheet argument is the source sheet itself).
# ---------------------------------------------------------------------------
$actual = $wb.Worksheets.Item("30_actual")
$actual.Copy([Type]::Missing, $actual)
$hard = $wb.Worksheets.Item("30_actual (2)")
$hard.Name = "30_actual_hard"

$genInfo = $wb.Worksheets.Item("GenInfo")

# ---------------------------------------------------------------------------
# 2. Data edits on the new "30_actual_hard" sheet + matching GenInfo notes.
#    Laura Jones (row 30) -> Laura James, and Alex Rodriguez's (row 31)
#    address is updated.
# ---------------------------------------------------------------------------
$hard.Range("A30").Value = "Laura James"
$genInfo.Range("A3").Value = "laura jones changed to laura james in 30 actual hard row 30"

$hard.Range("D31").Value = "8998 Dixon Rd, Las Colinas, IA 33229"
$genInfo.Range("A4").Value = "alex rodriguez address changed"

# ---------------------------------------------------------------------------
# 3. Restore each sheet's view/selection state.
# ---------------------------------------------------------------------------

# Best-effort: restore the workbook window geometry recorded in the saved
# file (harmless no-op if the host doesn't persist window placement).
try {
    $win = $wb.Windows.Item(1)
    $win.Left = 7340
    $win.Top = 880
    $win.Width = 20740
    $win.Height = 16320
} catch {}

$allCustomers = $wb.Worksheets.Item("All_customers")
$allCustomers.Activate()
try { $excel.ActiveWindow.ScrollRow = 11 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 1 } catch {}
$allCustomers.Range("A30").Select()

$counted24 = $wb.Worksheets.Item("24_counted")
$counted24.Activate()
$counted24.Range("C26").Select()

$actual.Activate()
$actual.Range("A30").Select()

$hard.Activate()
$hard.Range("H20").Select()

# "All customers (Table)" sheet view/selection is unchanged from the source
# workbook, so it is left untouched here.

# GenInfo becomes the active (last) tab, matching the saved workbook.
$genInfo.Activate()
$genInfo.Range("A5").Select()
